$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the original inlineStr cell type,
# so numeric-looking price strings are not silently reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.485.80'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").Value = '1.681.67'
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '216.75'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("D6").Value = '0.5326'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.2689'
$ws.Range("E8").Value = '  +4.00%  '
$ws.Range("D9").Value = '0.06403'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("D10").Value = '21.73'
$ws.Range("E10").Value = '  +5.76%  '
$ws.Range("D11").Value = '0.07795'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").Value = '1.682.05'
$ws.Range("E12").Value = '  +3.03%  '
$ws.Range("D13").Value = '4.496'
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").Value = '0.5578'
$ws.Range("E14").Value = '  +1.64%  '
$ws.Range("D15").Value = '0.0₅8327'
$ws.Range("E15").Value = '  +4.09%  '
$ws.Range("D16").Value = '65.66'
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = '26.532.11'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D19").Value = '4.765'
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = '194.42'
$ws.Range("E20").Value = '  +4.84%  '
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("D22").Value = '6.356'
$ws.Range("E22").Value = '  +4.02%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '143.12'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '0.1280'
$ws.Range("D26").Value = '7.442'
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("D27").Value = '16.32'
$ws.Range("E27").Value = '  +4.39%  '
$ws.Range("D28").Value = '1.426'
$ws.Range("E28").Value = '  +3.69%  '
$ws.Range("D29").Value = '0.06220'
$ws.Range("E29").Value = '  +5.03%  '
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("D31").Value = '3.607'
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("D32").Value = '3.453'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("D34").Value = '1.010'
$ws.Range("E34").Value = '  +2.98%  '
$ws.Range("D35").Value = '2.427'
$ws.Range("D36").Value = '2.791'
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").Value = '0.5741'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").Value = '0.01638'
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("D39").Value = '6.042'
$ws.Range("E39").Value = '  +6.59%  '
$ws.Range("D40").Value = '1.074.32'
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("D41").Value = '0.8578'
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").Value = '1.827.89'
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").Value = '57.20'
$ws.Range("E45").Value = '  +4.18%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.117'
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05207'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '6.029'
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.4239'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.461'
$ws.Range("E51").Value = '  +5.44%  '
